$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows with new dates
$ws.Range("B2").Value = "Jueves 16/05/2024"
$ws.Range("B3").Value = "Viernes 17/05/2024"
$ws.Range("B4").Value = "Lunes 20/05/2024"

# Add new rows 5-11
$data = @(
    @("2024", "Martes 21/05/2024"),
    @("2024", "Miércoles 22/05/2024"),
    @("2024", "Jueves 23/05/2024"),
    @("2024", "Lunes 27/05/2024"),
    @("2024", "Martes 28/05/2024"),
    @("2024", "Miércoles 29/05/2024"),
    @("2024", "Jueves 30/05/2024")
)

$row = 5
foreach ($entry in $data) {
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $entry[0]
    $cellA.Style = "Normal"
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}
